$d = $word.ActiveDocument

$replacements = @(
    @{old="2024-05-07 Tuesday"; new="2024-05-08 Wednesday"},
    @{old="72×36="; new="32×27="},
    @{old="50×76="; new="67×20="},
    @{old="31×26="; new="55×40="},
    @{old="23×15="; new="44×59="},
    @{old="58×21="; new="46×57="},
    @{old="16×28="; new="25×79="},
    @{old="46×44="; new="56×17="},
    @{old="57×85="; new="91×63="},
    @{old="87×22="; new="44×75="},
    @{old="24×67="; new="95×60="},
    @{old="37×44="; new="68×28="},
    @{old="55×95="; new="85×14="},
    @{old="56×83="; new="64×49="},
    @{old="85×41="; new="82×24="},
    @{old="17×55="; new="27×66="},
    @{old="44×53="; new="98×24="},
    @{old="17×78="; new="47×76="},
    @{old="67×13="; new="63×21="},
    @{old="96×44="; new="47×94="},
    @{old="72×14="; new="35×76="},
    @{old="61×38="; new="13×13="},
    @{old="52×15="; new="30×19="},
    @{old="28×36="; new="73×29="},
    @{old="11×33="; new="25×24="},
    @{old="81×14="; new="39×51="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
